$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, forcing text to preserve formatting (avoid numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.374.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0613"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0876"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.889.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.648.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.381.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.439.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.910"
$ws.Range("D37").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.797.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("E18").Value = "  -7.29%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -1.25%  "
